# Applies the cryptos.xlsx data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.485.37"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.33%  "

# Row 3: D3, E3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.765.85"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.86%  "

# Row 4: E4
$ws.Range("E4").Value = "  +0.79%  "

# Row 5: D5, E5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "598.20"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "

# Row 6: D6, E6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "162.50"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.80%  "

# Row 7: D7, E7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.762.12"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.96%  "

# Row 8: E8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9: D9, E9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.511"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "

# Row 10: D10, E10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.155"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.35%  "

# Row 11: D11, E11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.442"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "

# Row 12: D12, E12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "6.58"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +4.74%  "

# Row 13: D13, E13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000243"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.79%  "

# Row 14: D14, E14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "34.95"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.67%  "

# Row 15: D15, E15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.399.43"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.91%  "

# Row 16: D16, E16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.792.56"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17: D17, E17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "67.615.04"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

# Row 18: D18, E18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "18.15"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.98%  "

# Row 19: E19
$ws.Range("E19").Value = "  +1.77%  "

# Row 20: D20, E20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.96"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "

# Row 21: D21, E21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "455.64"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22: D22, E22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.43"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -4.64%  "

# Row 23: D23, E23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.688"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "

# Row 24: D24, E24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "82.55"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.85%  "

# Row 25: D25, E25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.0000141"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -6.49%  "

# Row 26: D26, E26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.79"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.17%  "

# Row 27: D27, E27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.07"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.09%  "

# Row 28: E28
$ws.Range("E28").Value = "  +0.01%  "

# Row 29: D29, E29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.79"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.09%  "

# Row 30: D30, E30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.916.76"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "

# Row 31: D31, E31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.79%  "

# Row 32: D32, E32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.91%  "

# Row 33: D33, E33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.57"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -6.79%  "

# Row 34: D34, E34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "28.75"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.50%  "

# Row 35: D35, E35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36: D36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "8.90"
$cell.Style = "Normal"

# Row 37: D37, E37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0985"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.46%  "

# Row 38: E38
$ws.Range("E38").Value = "  +3.59%  "

# Row 39: D39, E39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.74"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "

# Row 40: D40, E40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.975"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.49%  "

# Row 41: B41, C41, D41, E41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

# Row 42: B42, C42, D42, E42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.14"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -6.46%  "

# Row 44: D44, E44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "43.24"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.99%  "

# Row 45: D45, E45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "47.10"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.07%  "

# Row 46: D46, E46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "151.96"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.94%  "

# Row 47: D47, E47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.293"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.46%  "

# Row 48: B48, C48, D48, E48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.36"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49: B49, C49, D49, E49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.25"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50: D50, E50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.82"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.02%  "

# Row 51: D51, E51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "383.99"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.59%  "
